$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.94498
$ws.Range("H2").Value = 17.83494
$ws.Range("I2").Value = 0.4679240463447598
$ws.Range("J2").Value = 0.4679240463447597
$ws.Range("M2").Value = 0.428743
$ws.Range("N2").Value = 1.286229
$ws.Range("O2").Value = 0.00412050394863168
$ws.Range("P2").Value = 0.00412050394863168
$ws.Range("Q2").Value = 2.54886856014
$ws.Range("R2").Value = 22.93981704126
$ws.Range("S2").Value = 0.001928082880623296
$ws.Range("T2").Value = 0.001928082880623296
# Row 3
$ws.Range("G3").Value = 5.94498
$ws.Range("H3").Value = 17.83494
$ws.Range("I3").Value = 0.4679240463447598
$ws.Range("J3").Value = 0.4679240463447597
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("O3").Value = 0.7710272268990069
$ws.Range("P3").Value = 0.7710272268990069
$ws.Range("Q3").Value = 476.9433744402599
$ws.Range("R3").Value = 4292.490369962339
$ws.Range("S3").Value = 0.3607821798525625
$ws.Range("T3").Value = 0.3607821798525625
# Row 4
$ws.Range("G4").Value = 5.94498
$ws.Range("H4").Value = 17.83494
$ws.Range("I4").Value = 0.4679240463447598
$ws.Range("J4").Value = 0.4679240463447597
$ws.Range("M4").Value = 23.39612766666667
$ws.Range("N4").Value = 70.188383
$ws.Range("O4").Value = 0.2248522691523614
$ws.Range("P4").Value = 0.2248522691523614
$ws.Range("Q4").Value = 139.08951105578
$ws.Range("R4").Value = 1251.80559950202
$ws.Range("S4").Value = 0.105213783611574
$ws.Range("T4").Value = 0.1052137836115739
# Row 5
$ws.Range("H5").Value = 8.352077
$ws.Range("I5").Value = 0.219128164447035
$ws.Range("J5").Value = 0.219128164447035
$ws.Range("M5").Value = 0.428743
$ws.Range("N5").Value = 1.286229
$ws.Range("O5").Value = 0.00412050394863168
$ws.Range("P5").Value = 0.00412050394863168
$ws.Range("Q5").Value = 1.193631516403667
$ws.Range("R5").Value = 10.742683647633
$ws.Range("S5").Value = 0.0009029184668604197
$ws.Range("T5").Value = 0.0009029184668604197
# Row 6
$ws.Range("H6").Value = 8.352077
$ws.Range("I6").Value = 0.219128164447035
$ws.Range("J6").Value = 0.219128164447035
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("O6").Value = 0.7710272268990069
$ws.Range("P6").Value = 0.7710272268990069
$ws.Range("Q6").Value = 223.3519029480829
$ws.Range("S6").Value = 0.1689537809690669
$ws.Range("T6").Value = 0.1689537809690669
# Row 7
$ws.Range("H7").Value = 8.352077
$ws.Range("I7").Value = 0.219128164447035
$ws.Range("J7").Value = 0.219128164447035
$ws.Range("M7").Value = 23.39612766666667
$ws.Range("N7").Value = 70.188383
$ws.Range("O7").Value = 0.2248522691523614
$ws.Range("P7").Value = 0.2248522691523614
$ws.Range("Q7").Value = 65.13541992461012
$ws.Range("R7").Value = 586.218779321491
$ws.Range("S7").Value = 0.04927146501110762
$ws.Range("T7").Value = 0.04927146501110762
# Row 8
$ws.Range("G8").Value = 3.976005
$ws.Range("H8").Value = 11.928015
$ws.Range("I8").Value = 0.3129477892082053
$ws.Range("J8").Value = 0.3129477892082053
$ws.Range("M8").Value = 0.428743
$ws.Range("N8").Value = 1.286229
$ws.Range("O8").Value = 0.00412050394863168
$ws.Range("P8").Value = 0.00412050394863168
$ws.Range("Q8").Value = 1.704684311715
$ws.Range("R8").Value = 15.342158805435
$ws.Range("S8").Value = 0.001289502601147965
$ws.Range("T8").Value = 0.001289502601147965
# Row 9
$ws.Range("G9").Value = 3.976005
$ws.Range("H9").Value = 11.928015
$ws.Range("I9").Value = 0.3129477892082053
$ws.Range("J9").Value = 0.3129477892082053
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("O9").Value = 0.7710272268990069
$ws.Range("P9").Value = 0.7710272268990069
$ws.Range("Q9").Value = 318.9799194431849
$ws.Range("R9").Value = 2870.819274988665
$ws.Range("S9").Value = 0.2412912660773775
$ws.Range("T9").Value = 0.2412912660773775
# Row 10
$ws.Range("G10").Value = 3.976005
$ws.Range("H10").Value = 11.928015
$ws.Range("I10").Value = 0.3129477892082053
$ws.Range("J10").Value = 0.3129477892082053
$ws.Range("M10").Value = 23.39612766666667
$ws.Range("N10").Value = 70.188383
$ws.Range("O10").Value = 0.2248522691523614
$ws.Range("P10").Value = 0.2248522691523614
$ws.Range("Q10").Value = 93.02312058330502
$ws.Range("R10").Value = 837.2080852497451
$ws.Range("S10").Value = 0.07036702052967984
$ws.Range("T10").Value = 0.07036702052967984
